$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Capirava" -> "Capivara" in cell A16
$ws.Range("A16").Value = "Capivara"

# Update the last active selection to match the saved view state
$ws.Range("E17:E18").Select()
